$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 (I0) and J1 (IF), copying the style of the existing
# header cell H1 so the new headers match (bold, centered, bordered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for column I (I0) and column J (IF), rows 2-56.
$iValues = @(6,8,9,9,1,6,8,7,10,7,8,8,8,9,8,7,6,11,8,8,8,6,7,8,4,8,7,7,7,7,5,8,5,7,10,8,6,7,7,7,8,8,8,7,6,5,6,6,6,8,4,8,7,6,7)
$jValues = @(6,8,9,9,1,6,8,7,10,8,8,9,8,9,8,7,7,11,8,8,8,6,8,8,4,8,8,7,7,8,6,9,6,7,10,8,6,8,7,7,8,8,8,7,6,6,6,6,6,8,4,8,7,6,7)

for ($i = 0; $i -lt $iValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $iValues[$i]
    $ws.Cells.Item($row, 10).Value = $jValues[$i]
}
